$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-31 in columns A and B contain a linearly spaced sequence of 30 values.
# Previously it spanned 1 .. 259; now it should span 100 .. 150.
$start = 100.0
$end = 150.0
$n = 30

for ($i = 0; $i -lt $n; $i++) {
    $value = $start + ($end - $start) * $i / ($n - 1)
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $value
    $ws.Cells.Item($row, 2).Value = $value
}
